$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows after row 22 to hold the data that is being
#     displaced from rows 21/22 (Segunda / Tercera, previous week) ---
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(23).Insert()

# Row 23 (new): copy of old row 21 (Segunda, previous week 44235)
$ws.Cells.Item(23,1).Value = 3
$ws.Cells.Item(23,2).Value = "Femacal de La Calera"
$ws.Cells.Item(23,3).Value = "Coquimbo"
$ws.Cells.Item(23,4).Value = 44235
$ws.Cells.Item(23,5).Value = 5
$ws.Cells.Item(23,6).Value = 100112043
$ws.Cells.Item(23,7).Value = "Pepino dulce"
$ws.Cells.Item(23,8).Value = "Cultivar IV Región"
$ws.Cells.Item(23,9).Value = "Segunda"
$ws.Cells.Item(23,10).Value = 70
$ws.Cells.Item(23,11).Value = 12000
$ws.Cells.Item(23,12).Value = 12000
$ws.Cells.Item(23,13).Value = 12000
$ws.Cells.Item(23,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(23,15).Value = "Provincia de Limarí"
$ws.Cells.Item(23,16).Value = 667
$ws.Cells.Item(23,17).Value = 18
$ws.Cells.Item(23,18).Value = "Hortaliza"
$ws.Cells.Item(23,4).NumberFormat = $ws.Cells.Item(20,4).NumberFormat

# Row 24 (new): copy of old row 22 (Tercera, previous week 44235)
$ws.Cells.Item(24,1).Value = 3
$ws.Cells.Item(24,2).Value = "Femacal de La Calera"
$ws.Cells.Item(24,3).Value = "Coquimbo"
$ws.Cells.Item(24,4).Value = 44235
$ws.Cells.Item(24,5).Value = 5
$ws.Cells.Item(24,6).Value = 100112043
$ws.Cells.Item(24,7).Value = "Pepino dulce"
$ws.Cells.Item(24,8).Value = "Cultivar IV Región"
$ws.Cells.Item(24,9).Value = "Tercera"
$ws.Cells.Item(24,10).Value = 60
$ws.Cells.Item(24,11).Value = 10000
$ws.Cells.Item(24,12).Value = 10000
$ws.Cells.Item(24,13).Value = 10000
$ws.Cells.Item(24,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(24,15).Value = "Provincia de Limarí"
$ws.Cells.Item(24,16).Value = 556
$ws.Cells.Item(24,17).Value = 18
$ws.Cells.Item(24,18).Value = "Hortaliza"
$ws.Cells.Item(24,4).NumberFormat = $ws.Cells.Item(20,4).NumberFormat

# --- Row 20 (Primera): update to new week's figures ---
$ws.Cells.Item(20,4).Value = 45132
$ws.Cells.Item(20,8).Value = "Sin especificar"
$ws.Cells.Item(20,10).Value = 50
$ws.Cells.Item(20,11).Value = 18000
$ws.Cells.Item(20,12).Value = 18000
$ws.Cells.Item(20,13).Value = 18000
$ws.Cells.Item(20,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(20,16).Value = 1200
$ws.Cells.Item(20,17).Value = 15

# --- Row 21 (Segunda): update to new week's figures ---
$ws.Cells.Item(21,4).Value = 45132
$ws.Cells.Item(21,8).Value = "Sin especificar"
$ws.Cells.Item(21,10).Value = 45
$ws.Cells.Item(21,11).Value = 13000
$ws.Cells.Item(21,12).Value = 13000
$ws.Cells.Item(21,13).Value = 13000
$ws.Cells.Item(21,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(21,16).Value = 867
$ws.Cells.Item(21,17).Value = 15

# --- Row 22 (was Tercera, becomes the new week's Primera-grade entry) ---
$ws.Cells.Item(22,9).Value = "Primera"
$ws.Cells.Item(22,10).Value = 80
$ws.Cells.Item(22,11).Value = 14000
$ws.Cells.Item(22,12).Value = 14000
$ws.Cells.Item(22,13).Value = 14000
$ws.Cells.Item(22,16).Value = 778

$wb.Save()
